$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these cells remain text (matching the source data which stores
# prices as inline strings), not auto-converted to numbers by Excel.
$ws.Range("D2:D50").NumberFormat = "@"

$ws.Range("D2").Value = "271.93"
$ws.Range("D3").Value = "23.13"
$ws.Range("D4").Value = "6.380"
$ws.Range("D5").Value = "0.06283"
$ws.Range("D6").Value = "3.659"
$ws.Range("D7").Value = "6.727"
$ws.Range("D8").Value = "1.386"
$ws.Range("D10").Value = "0.1637"
$ws.Range("D11").Value = "0.08443"
$ws.Range("D12").Value = "0.03490"
$ws.Range("D13").Value = "0.03138"
$ws.Range("D14").Value = "0.09319"
$ws.Range("D15").Value = "3.877"
$ws.Range("D16").Value = "0.001729"
$ws.Range("D17").Value = "0.04832"
$ws.Range("D18").Value = "0.006265"
$ws.Range("D19").Value = "0.005351"
$ws.Range("D20").Value = "0.001088"
$ws.Range("D22").Value = "3.733"
$ws.Range("D23").Value = "2.320"
$ws.Range("D24").Value = "0.01381"
$ws.Range("D25").Value = "0.3407"
$ws.Range("D26").Value = "0.1263"
$ws.Range("D40").Value = "0.04691"
$ws.Range("D41").Value = "0.006906"
$ws.Range("D45").Value = "0.00006243"
$ws.Range("D46").Value = "0.00000000748"
$ws.Range("D47").Value = "0.7979"
$ws.Range("D48").Value = "0.09527"
$ws.Range("D50").Value = "0.01237"

# Restore default (Normal) style so no stray number-format style lingers
# on the edited cells.
$ws.Range("D2:D50").Style = "Normal"

